$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new test case row (pulled in from JSON test data).
# Set column B first so the shared-string table records "Y" before
# "LoginTest", matching how Excel lays out the strings for this edit.
$ws.Range("B4").Value = "Y"
$ws.Range("A4").Value = "LoginTest"

# Move the selection like Excel would after navigating away from the edit
$ws.Range("A8").Select()
